$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New resume rows appended below the existing data (rows 53-57).
$newRows = @(
    @{ A = "SANTOSH  KUMAR  DASH"; B = "Santoshkdash27@gmail.com"; C = "2024"; D = "1 years";       E = 38.9;  F = "+91 8197576872" },
    @{ A = "Soumya Ranjan Swain";  B = "soumya18.swain@gmail.com"; C = "2023"; D = "No Experience";  E = 43.02; F = "+91-7077964867" },
    @{ A = "SANTOSH  KUMAR  DASH"; B = "Santoshkdash27@gmail.com"; C = "2024"; D = "1 years";       E = 39.07; F = "+91 8197576872" },
    @{ A = "Soumya Ranjan Swain";  B = "soumya18.swain@gmail.com"; C = "2023"; D = "No Experience";  E = 43.02; F = "+91-7077964867" },
    @{ A = "Job Description";     B = "Not Found";               C = "Not Found"; D = "No Experience"; E = 100;   F = "Not Found" }
)

$startRow = 53
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($row, 1).Value = $data.A
    $ws.Cells.Item($row, 2).Value = $data.B

    # Keep Batch Year / Job Description text columns as text (not auto-converted to numbers).
    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = $data.C

    $ws.Cells.Item($row, 4).Value = $data.D
    $ws.Cells.Item($row, 5).Value = $data.E
    $ws.Cells.Item($row, 6).Value = $data.F
}
